$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.441.42'
$ws.Range("E2").Value = '  -2.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.271.89'
$ws.Range("E3").Value = '  -4.96%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.02'
$ws.Range("E5").Value = '  -5.51%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '605.44'
$ws.Range("E6").Value = '  -5.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.35'
$ws.Range("E7").Value = '  -7.70%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.374'
$ws.Range("E8").Value = '  -6.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.904'
$ws.Range("E10").Value = '  -8.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.263.13'
$ws.Range("E11").Value = '  -5.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '41.41'
$ws.Range("E12").Value = '  +0.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.190'
$ws.Range("E13").Value = '  -3.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '92.244.52'
$ws.Range("E14").Value = '  -2.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.86'
$ws.Range("E15").Value = '  -3.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.881.25'
$ws.Range("E16").Value = '  -5.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000238'
$ws.Range("E17").Value = '  -6.77%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.87'
$ws.Range("E18").Value = '  -6.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.267.62'
$ws.Range("E19").Value = '  -4.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.90'
$ws.Range("E20").Value = '  -4.46%  '

$ws.Range("E21").Value = '  -7.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.38'
$ws.Range("E22").Value = '  +7.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '482.10'
$ws.Range("E23").Value = '  -3.69%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.433'
$ws.Range("E24").Value = '  -14.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000175'
$ws.Range("E25").Value = '  -7.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.95'
$ws.Range("E26").Value = '  -9.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '88.34'
$ws.Range("E27").Value = '  -3.75%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.48'
$ws.Range("E28").Value = '  -3.86%  '

$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.440.80'
$ws.Range("E29").Value = '  -5.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.82'
$ws.Range("E31").Value = '  -7.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.136'
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.57'
$ws.Range("E34").Value = '  -6.27%  '

$ws.Range("E35").Value = '  -7.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '27.64'
$ws.Range("E36").Value = '  -10.80%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.516'
$ws.Range("E37").Value = '  -8.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '528.36'
$ws.Range("E38").Value = '  +1.07%  '

$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.17'
$ws.Range("E40").Value = '  -6.32%  '

$ws.Range("E41").Value = '  -3.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.33'
$ws.Range("E42").Value = '  -6.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.90'
$ws.Range("E43").Value = '  -0.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.839'
$ws.Range("E44").Value = '  -7.47%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.55'
$ws.Range("E45").Value = '  +2.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.64'
$ws.Range("E46").Value = '  -3.14%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0402'
$ws.Range("E47").Value = '  -2.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.21'
$ws.Range("E48").Value = '  -6.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '51.38'
$ws.Range("E49").Value = '  -3.96%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.05'
$ws.Range("E50").Value = '  -4.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.70'
$ws.Range("E51").Value = '  -3.50%  '
